$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 104593637
$ws.Range("B2").Value = 96334
$ws.Range("D2").Value = "VU"
$ws.Range("E2").Value = 220787
$ws.Range("F2").Value = "Knärot"
$ws.Range("G2").Value = "Goodyera repens"
$ws.Range("H2").Value = "(L.) R. Br."
$ws.Range("K2").ClearContents() | Out-Null
$ws.Range("L2").ClearContents() | Out-Null
$ws.Range("M2").ClearContents() | Out-Null
$ws.Range("N2").ClearContents() | Out-Null
$ws.Range("Q2").Value = 473780.3094452888
$ws.Range("R2").Value = 7013777.427834951
$ws.Range("AC2").ClearContents() | Out-Null

# Row 3
$ws.Range("A3").Value = 104593632
$ws.Range("B3").Value = 96334
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 220787
$ws.Range("F3").Value = "Knärot"
$ws.Range("G3").Value = "Goodyera repens"
$ws.Range("H3").Value = "(L.) R. Br."
$ws.Range("K3").ClearContents() | Out-Null
$ws.Range("L3").ClearContents() | Out-Null
$ws.Range("M3").ClearContents() | Out-Null
$ws.Range("N3").ClearContents() | Out-Null
$ws.Range("Q3").Value = 473815.7661137963
$ws.Range("R3").Value = 7013977.153526685
$ws.Range("AC3").ClearContents() | Out-Null

# Row 4
$ws.Range("A4").Value = 104593622
$ws.Range("B4").Value = 56395
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("K4").Value = "'"
$ws.Range("K4").Style = "Normal"
$ws.Range("L4").Value = "'"
$ws.Range("L4").Style = "Normal"
$ws.Range("M4").Value = "äldre spår"
$ws.Range("N4").Value = "'"
$ws.Range("N4").Style = "Normal"
$ws.Range("Q4").Value = 473830.905343441
$ws.Range("R4").Value = 7013897.29666794
$ws.Range("AC4").Value = "ringhack"

# Row 5
$ws.Range("A5").Value = 104593631
$ws.Range("B5").Value = 96334
$ws.Range("D5").Value = "VU"
$ws.Range("E5").Value = 220787
$ws.Range("F5").Value = "Knärot"
$ws.Range("G5").Value = "Goodyera repens"
$ws.Range("H5").Value = "(L.) R. Br."
$ws.Range("K5").ClearContents() | Out-Null
$ws.Range("L5").ClearContents() | Out-Null
$ws.Range("M5").ClearContents() | Out-Null
$ws.Range("N5").ClearContents() | Out-Null
$ws.Range("Q5").Value = 473812.0075608135
$ws.Range("R5").Value = 7013958.714830574
$ws.Range("AC5").ClearContents() | Out-Null

# Row 6
$ws.Range("A6").Value = 104593630
$ws.Range("B6").Value = 96334
$ws.Range("D6").Value = "VU"
$ws.Range("E6").Value = 220787
$ws.Range("F6").Value = "Knärot"
$ws.Range("G6").Value = "Goodyera repens"
$ws.Range("H6").Value = "(L.) R. Br."
$ws.Range("K6").ClearContents() | Out-Null
$ws.Range("L6").ClearContents() | Out-Null
$ws.Range("M6").ClearContents() | Out-Null
$ws.Range("N6").ClearContents() | Out-Null
$ws.Range("Q6").Value = 473798.8866438381
$ws.Range("R6").Value = 7013953.866335354
$ws.Range("AC6").ClearContents() | Out-Null

# Row 7
$ws.Range("A7").Value = 104593624
$ws.Range("B7").Value = 96334
$ws.Range("D7").Value = "VU"
$ws.Range("E7").Value = 220787
$ws.Range("F7").Value = "Knärot"
$ws.Range("G7").Value = "Goodyera repens"
$ws.Range("H7").Value = "(L.) R. Br."
$ws.Range("K7").ClearContents() | Out-Null
$ws.Range("L7").ClearContents() | Out-Null
$ws.Range("M7").ClearContents() | Out-Null
$ws.Range("N7").ClearContents() | Out-Null
$ws.Range("Q7").Value = 473801.0947980214
$ws.Range("R7").Value = 7013892.583679659
$ws.Range("AC7").ClearContents() | Out-Null

# Row 8
$ws.Range("A8").Value = 104593620
$ws.Range("B8").Value = 89392
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 1202
$ws.Range("F8").Value = "Ullticka"
$ws.Range("G8").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H8").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("K8").ClearContents() | Out-Null
$ws.Range("L8").ClearContents() | Out-Null
$ws.Range("M8").ClearContents() | Out-Null
$ws.Range("N8").ClearContents() | Out-Null
$ws.Range("Q8").Value = 473808.0725733605
$ws.Range("R8").Value = 7013974.062789564
$ws.Range("AC8").ClearContents() | Out-Null

# Row 9
$ws.Range("A9").Value = 104593621
$ws.Range("B9").Value = 89392
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 1202
$ws.Range("F9").Value = "Ullticka"
$ws.Range("G9").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H9").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("K9").ClearContents() | Out-Null
$ws.Range("L9").ClearContents() | Out-Null
$ws.Range("M9").ClearContents() | Out-Null
$ws.Range("N9").ClearContents() | Out-Null
$ws.Range("Q9").Value = 473766.1571259646
$ws.Range("R9").Value = 7013701.408301079
$ws.Range("AC9").ClearContents() | Out-Null

# Row 10
$ws.Range("A10").Value = 104593623
$ws.Range("B10").Value = 56395
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 100109
$ws.Range("F10").Value = "Tretåig hackspett"
$ws.Range("G10").Value = "Picoides tridactylus"
$ws.Range("H10").Value = "(Linnaeus, 1758)"
$ws.Range("K10").Value = "'"
$ws.Range("K10").Style = "Normal"
$ws.Range("L10").Value = "'"
$ws.Range("L10").Style = "Normal"
$ws.Range("M10").Value = "äldre spår"
$ws.Range("N10").Value = "'"
$ws.Range("N10").Style = "Normal"
$ws.Range("Q10").Value = 473722.3678416939
$ws.Range("R10").Value = 7013918.902128431
$ws.Range("AC10").Value = "ringhack"

# Row 11
$ws.Range("A11").Value = 104593627
$ws.Range("B11").Value = 96334
$ws.Range("D11").Value = "VU"
$ws.Range("E11").Value = 220787
$ws.Range("F11").Value = "Knärot"
$ws.Range("G11").Value = "Goodyera repens"
$ws.Range("H11").Value = "(L.) R. Br."
$ws.Range("K11").ClearContents() | Out-Null
$ws.Range("L11").ClearContents() | Out-Null
$ws.Range("M11").ClearContents() | Out-Null
$ws.Range("N11").ClearContents() | Out-Null
$ws.Range("Q11").Value = 473701.5160585373
$ws.Range("R11").Value = 7013906.458910029
$ws.Range("AC11").ClearContents() | Out-Null

# Row 12
$ws.Range("A12").Value = 104593628
$ws.Range("B12").Value = 96334
$ws.Range("D12").Value = "VU"
$ws.Range("E12").Value = 220787
$ws.Range("F12").Value = "Knärot"
$ws.Range("G12").Value = "Goodyera repens"
$ws.Range("H12").Value = "(L.) R. Br."
$ws.Range("K12").ClearContents() | Out-Null
$ws.Range("L12").ClearContents() | Out-Null
$ws.Range("M12").ClearContents() | Out-Null
$ws.Range("N12").ClearContents() | Out-Null
$ws.Range("Q12").Value = 473726.4678040863
$ws.Range("R12").Value = 7013923.823926651
$ws.Range("AC12").ClearContents() | Out-Null

# Row 13
$ws.Range("A13").Value = 104593626
$ws.Range("B13").Value = 96334
$ws.Range("D13").Value = "VU"
$ws.Range("E13").Value = 220787
$ws.Range("F13").Value = "Knärot"
$ws.Range("G13").Value = "Goodyera repens"
$ws.Range("H13").Value = "(L.) R. Br."
$ws.Range("K13").ClearContents() | Out-Null
$ws.Range("L13").ClearContents() | Out-Null
$ws.Range("M13").ClearContents() | Out-Null
$ws.Range("N13").ClearContents() | Out-Null
$ws.Range("Q13").Value = 473718.6013391476
$ws.Range("R13").Value = 7013899.562304306
$ws.Range("AC13").ClearContents() | Out-Null

# Row 14
$ws.Range("A14").Value = 104593634
$ws.Range("B14").Value = 96334
$ws.Range("D14").Value = "VU"
$ws.Range("E14").Value = 220787
$ws.Range("F14").Value = "Knärot"
$ws.Range("G14").Value = "Goodyera repens"
$ws.Range("H14").Value = "(L.) R. Br."
$ws.Range("K14").ClearContents() | Out-Null
$ws.Range("L14").ClearContents() | Out-Null
$ws.Range("M14").ClearContents() | Out-Null
$ws.Range("N14").ClearContents() | Out-Null
$ws.Range("Q14").Value = 473769.3477768434
$ws.Range("R14").Value = 7013705.43688098
$ws.Range("AC14").ClearContents() | Out-Null

# Row 15
$ws.Range("A15").Value = 104593636
$ws.Range("B15").Value = 96334
$ws.Range("D15").Value = "VU"
$ws.Range("E15").Value = 220787
$ws.Range("F15").Value = "Knärot"
$ws.Range("G15").Value = "Goodyera repens"
$ws.Range("H15").Value = "(L.) R. Br."
$ws.Range("K15").ClearContents() | Out-Null
$ws.Range("L15").ClearContents() | Out-Null
$ws.Range("M15").ClearContents() | Out-Null
$ws.Range("N15").ClearContents() | Out-Null
$ws.Range("Q15").Value = 473782.4035598941
$ws.Range("R15").Value = 7013757.588904253
$ws.Range("AC15").ClearContents() | Out-Null

# Row 16
$ws.Range("A16").Value = 104593625
$ws.Range("B16").Value = 96334
$ws.Range("D16").Value = "VU"
$ws.Range("E16").Value = 220787
$ws.Range("F16").Value = "Knärot"
$ws.Range("G16").Value = "Goodyera repens"
$ws.Range("H16").Value = "(L.) R. Br."
$ws.Range("K16").ClearContents() | Out-Null
$ws.Range("L16").ClearContents() | Out-Null
$ws.Range("M16").ClearContents() | Out-Null
$ws.Range("N16").ClearContents() | Out-Null
$ws.Range("Q16").Value = 473775.8828205758
$ws.Range("R16").Value = 7013898.645042086
$ws.Range("AC16").ClearContents() | Out-Null

# Row 17
$ws.Range("A17").Value = 104593629
$ws.Range("B17").Value = 96334
$ws.Range("D17").Value = "VU"
$ws.Range("E17").Value = 220787
$ws.Range("F17").Value = "Knärot"
$ws.Range("G17").Value = "Goodyera repens"
$ws.Range("H17").Value = "(L.) R. Br."
$ws.Range("K17").ClearContents() | Out-Null
$ws.Range("L17").ClearContents() | Out-Null
$ws.Range("M17").ClearContents() | Out-Null
$ws.Range("N17").ClearContents() | Out-Null
$ws.Range("Q17").Value = 473760.983223469
$ws.Range("R17").Value = 7013952.372943264
$ws.Range("AC17").ClearContents() | Out-Null
